$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.430.58"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.056.62"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.46"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.16"
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.00"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "2.362.43"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.91"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.758"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.31"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "2.055.15"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "37.279.43"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.76"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "0.0₃0829"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.93"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.86"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.129"
$ws.Range("E28").Value = "  -7.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.14"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.56"
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.63"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.28"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("D41").Value = "1.503.16"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.16"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.20"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0948"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.14"
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").Value = "2.247.11"
$ws.Range("E51").Value = "  -0.91%  "
